$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '311.44'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '-1.09%'
$c.Style = "Normal"
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '37.47'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '-4.58%'
$c.Style = "Normal"
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '5.067'
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '-1.59%'
$c.Style = "Normal"
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.07752'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '-5.13%'
$c.Style = "Normal"
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '4.362'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '-0.51%'
$c.Style = "Normal"
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("B7")
$c.NumberFormat = "@"
$c.Value = 'KuCoinToken'
$c.Style = "Normal"
$c = $ws.Range("C7")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '8.191'
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '-1.61%'
$c.Style = "Normal"
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("B8")
$c.NumberFormat = "@"
$c.Value = 'FTXToken'
$c.Style = "Normal"
$c = $ws.Range("C8")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '1.870'
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '-5.75%'
$c.Style = "Normal"
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '-9.17%'
$c.Style = "Normal"
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.9179'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '-1.89%'
$c.Style = "Normal"
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.1206'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '-7.05%'
$c.Style = "Normal"
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.1903'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '-3.95%'
$c.Style = "Normal"
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.08842'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '-1.76%'
$c.Style = "Normal"
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.03386'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '-4.00%'
$c.Style = "Normal"
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.09704'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '-0.34%'
$c.Style = "Normal"
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.001370'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '-2.28%'
$c.Style = "Normal"
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.006067'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '-8.24%'
$c.Style = "Normal"
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.543'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '-2.23%'
$c.Style = "Normal"
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.3405'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '-1.81%'
$c.Style = "Normal"
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.1281'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '-2.04%'
$c.Style = "Normal"
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '5.045'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '0.00%'
$c.Style = "Normal"
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.2596'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '4.31%'
$c.Style = "Normal"
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.02108'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '5,599.88%'
$c.Style = "Normal"
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.04394'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '0.30%'
$c.Style = "Normal"
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '-2.35%'
$c.Style = "Normal"
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.004238'
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '-10.88%'
$c.Style = "Normal"
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.0001353'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '-65.25%'
$c.Style = "Normal"
$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.02095'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '-6.80%'
$c.Style = "Normal"
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.04927'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '-5.49%'
$c.Style = "Normal"
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.007915'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '2.11%'
$c.Style = "Normal"
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.009943'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '-3.75%'
$c.Style = "Normal"
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.1340'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '-4.17%'
$c.Style = "Normal"
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.002064'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '-1.78%'
$c.Style = "Normal"
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.009654'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '8.72%'
$c.Style = "Normal"
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.00006587'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '-3.44%'
$c.Style = "Normal"
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.00000000752'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '0.13%'
$c.Style = "Normal"
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.003047'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '1.28%'
$c.Style = "Normal"
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '-0.15%'
$c.Style = "Normal"
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.00002104'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '0.13%'
$c.Style = "Normal"
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0002004'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '0.13%'
$c.Style = "Normal"
$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
